$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 165, shifting rows 165:174 down to 166:175.
$ws.Rows(165).Insert()

# Populate the newly inserted row 165 with the new weekly price record.
$ws.Cells.Item(165, 1).Value  = 4
$ws.Cells.Item(165, 2).Value  = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(165, 3).Value  = 'Los Lagos'
$ws.Cells.Item(165, 4).Value  = 44516
$ws.Cells.Item(165, 5).Value  = 10
$ws.Cells.Item(165, 6).Value  = 'Fruta'
$ws.Cells.Item(165, 7).Value  = 100102
$ws.Cells.Item(165, 8).Value  = 'Cítricos'
$ws.Cells.Item(165, 9).Value  = 100102006
$ws.Cells.Item(165, 10).Value = 'Pomelo'
$ws.Cells.Item(165, 11).Value = 'Start Ruby'
$ws.Cells.Item(165, 12).Value = 'Primera'
$ws.Cells.Item(165, 13).Value = 300
$ws.Cells.Item(165, 14).Value = 11000
$ws.Cells.Item(165, 15).Value = 12000
$ws.Cells.Item(165, 16).Value = 11500
$ws.Cells.Item(165, 17).Value = '$/caja 14 kilos empedrada'
$ws.Cells.Item(165, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(165, 19).Value = 821
$ws.Cells.Item(165, 20).Value = 14
